# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to various Leve profit-tracking sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 1482.0834
$ws.Range("I43").Value = 890
$ws.Range("K43").Value = 890
$ws.Range("M43").Value = -821

# Row 62
$ws.Range("H62").Value = 3666
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 3999
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 3999
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -5247

# Row 65
$ws.Range("H65").Value = 3666
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 3999
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 19995
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -26235

# Row 92
$ws.Range("H92").Value = 1368252.5
$ws.Range("I92").Value = 2052181.4
$ws.Range("J92").Value = 394.66666
$ws.Range("K92").Value = 2052181.4
$ws.Range("L92").Value = 394.66666
$ws.Range("M92").Value = -2050933.4
$ws.Range("N92").Value = -2890.66666

# Row 96
$ws.Range("H96").Value = 2763
$ws.Range("I96").Value = 2763
$ws.Range("K96").Value = 8289
$ws.Range("M96").Value = -6916

# Row 112
$ws.Range("H112").Value = 4222.1875
$ws.Range("I112").Value = 933.3333
$ws.Range("J112").Value = 4981.154
$ws.Range("K112").Value = 2799.9999
$ws.Range("L112").Value = 14943.462
$ws.Range("M112").Value = -1691.9999
$ws.Range("N112").Value = -17159.462

# Row 137
$ws.Range("H137").Value = 1484.2424
$ws.Range("I137").Value = 928
$ws.Range("J137").Value = 2151.7334
$ws.Range("K137").Value = 2784
$ws.Range("L137").Value = 6455.2002
$ws.Range("M137").Value = -234
$ws.Range("N137").Value = -11555.2002

# Row 138
$ws.Range("H138").Value = 2948.0833
$ws.Range("I138").Value = 3087.7
$ws.Range("K138").Value = 9263.099999999999
$ws.Range("M138").Value = -4123.099999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3324.261
$ws.Range("I32").Value = 2232.7307
$ws.Range("J32").Value = 6663.0586
$ws.Range("K32").Value = 2232.7307
$ws.Range("L32").Value = 6663.0586
$ws.Range("M32").Value = -1945.7307
$ws.Range("N32").Value = -7237.0586

# Row 45
$ws.Range("H45").Value = 1522.7
$ws.Range("I45").Value = 968.25
$ws.Range("K45").Value = 968.25
$ws.Range("M45").Value = -591.25

# Row 61
$ws.Range("H61").Value = 2734.9033
$ws.Range("I61").Value = 1788.8
$ws.Range("J61").Value = 6677
$ws.Range("K61").Value = 1788.8
$ws.Range("L61").Value = 6677
$ws.Range("M61").Value = -1576.8
$ws.Range("N61").Value = -7101

# Row 74
$ws.Range("H74").Value = 772.55554
$ws.Range("I74").Value = 771
$ws.Range("J74").Value = 775.6667
$ws.Range("K74").Value = 771
$ws.Range("L74").Value = 775.6667
$ws.Range("M74").Value = 103
$ws.Range("N74").Value = -2523.6667

# Row 77
$ws.Range("H77").Value = 772.55554
$ws.Range("I77").Value = 771
$ws.Range("J77").Value = 775.6667
$ws.Range("K77").Value = 3855
$ws.Range("L77").Value = 3878.3335
$ws.Range("M77").Value = 513
$ws.Range("N77").Value = -12614.3335

# Row 122
$ws.Range("H122").Value = 1700.65
$ws.Range("I122").Value = 1518.5454
$ws.Range("J122").Value = 1923.2222
$ws.Range("K122").Value = 4555.6362
$ws.Range("L122").Value = 5769.6666
$ws.Range("M122").Value = -2105.6362
$ws.Range("N122").Value = -10669.6666

# Row 136
$ws.Range("H136").Value = 2734.9033
$ws.Range("I136").Value = 1788.8
$ws.Range("J136").Value = 6677
$ws.Range("K136").Value = 5366.4
$ws.Range("L136").Value = 20031
$ws.Range("M136").Value = -2816.4
$ws.Range("N136").Value = -25131

$ws = $wb.Worksheets.Item("BSM")
# Row 75
$ws.Range("H75").Value = 1000
$ws.Range("I75").Value = 1000
$ws.Range("K75").Value = 1000
$ws.Range("M75").Value = -64

# Row 78
$ws.Range("H78").Value = 1000
$ws.Range("I78").Value = 1000
$ws.Range("K78").Value = 3000
$ws.Range("M78").Value = 1680

$ws = $wb.Worksheets.Item("CRP")
# Row 94
$ws.Range("H94").Value = 934.3077
$ws.Range("J94").Value = 1011.7143
$ws.Range("L94").Value = 1011.7143
$ws.Range("N94").Value = -1913.7143

# Row 134
$ws.Range("H134").Value = 903.1539
$ws.Range("I134").Value = 811.75
$ws.Range("K134").Value = 2435.25
$ws.Range("M134").Value = 99.75

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 765.13
$ws.Range("J131").Value = 807.9340999999999
$ws.Range("L131").Value = 2423.8023
$ws.Range("N131").Value = -12503.8023

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 2499.875
$ws.Range("I97").Value = 2225
$ws.Range("J97").Value = 2774.75
$ws.Range("K97").Value = 2225
$ws.Range("L97").Value = 2774.75
$ws.Range("M97").Value = -1729
$ws.Range("N97").Value = -3766.75

# Row 122
$ws.Range("H122").Value = 2385.4211
$ws.Range("I122").Value = 2268.4443
$ws.Range("J122").Value = 2490.7
$ws.Range("K122").Value = 6805.3329
$ws.Range("L122").Value = 7472.099999999999
$ws.Range("M122").Value = -4355.3329
$ws.Range("N122").Value = -12372.1

# Row 132
$ws.Range("H132").Value = 5339.2607
$ws.Range("I132").Value = 4291.8237
$ws.Range("J132").Value = 8307
$ws.Range("K132").Value = 12875.4711
$ws.Range("L132").Value = 24921
$ws.Range("M132").Value = -10345.4711
$ws.Range("N132").Value = -29981

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4290.316
$ws.Range("I7").Value = 1964.3636
$ws.Range("K7").Value = 1964.3636
$ws.Range("M7").Value = -1852.3636

# Row 55
$ws.Range("H55").Value = 528.7143
$ws.Range("J55").Value = 548.4
$ws.Range("L55").Value = 548.4
$ws.Range("N55").Value = -894.4

# Row 122
$ws.Range("H122").Value = 4339.4736
$ws.Range("I122").Value = 3055.5
$ws.Range("J122").Value = 5766.1113
$ws.Range("K122").Value = 9166.5
$ws.Range("L122").Value = 17298.3339
$ws.Range("M122").Value = -6716.5
$ws.Range("N122").Value = -22198.3339

# Row 126
$ws.Range("H126").Value = 4290.316
$ws.Range("I126").Value = 1964.3636
$ws.Range("K126").Value = 5893.0908
$ws.Range("M126").Value = -3423.0908

# Row 132
$ws.Range("H132").Value = 2141.1924
$ws.Range("I132").Value = 1436.1428
$ws.Range("J132").Value = 2400.9473
$ws.Range("K132").Value = 4308.428400000001
$ws.Range("L132").Value = 7202.841899999999
$ws.Range("M132").Value = -1778.428400000001
$ws.Range("N132").Value = -12262.8419

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 8825.913
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 8825.913
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 26477.739
$ws.Range("N132").Value = -31537.739
$ws.Range("M132").ClearContents()

# Row 133
$ws.Range("H133").Value = 55715
$ws.Range("J133").Value = 55715
$ws.Range("L133").Value = 55715
$ws.Range("N133").Value = -65835

# Row 136
$ws.Range("H136").Value = 2196.4375
$ws.Range("I136").Value = 1157.8334
$ws.Range("K136").Value = 3473.5002
$ws.Range("M136").Value = -923.5001999999999
